$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ResultsChart")

$ws.Range("A2").Value = 584146.0
$ws.Range("B2").Value = 7388145.0
$ws.Range("C2").Value = 1353054.0
$ws.Range("D2").Value = 573276.0
$ws.Range("E2").Value = 5422733.0
